# Auto-generated edit script applying numeric corrections described in the commit diff
# for workbook "Cactuar_Profits" (8 worksheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1257
$ws.Range("I41").Value = 662.75
$ws.Range("J41").Value = 1494.7
$ws.Range("K41").Value = 662.75
$ws.Range("L41").Value = 1494.7
$ws.Range("M41").Value = -222.75
$ws.Range("N41").Value = -2374.7
$ws.Range("H64").Value = 1504013
$ws.Range("I64").Value = 3110112.5
$ws.Range("K64").Value = 3110112.5
$ws.Range("M64").Value = -3109864.5
$ws.Range("H67").Value = 1504013
$ws.Range("I67").Value = 3110112.5
$ws.Range("K67").Value = 3110112.5
$ws.Range("M67").Value = -3109254.5
$ws.Range("H76").Value = 3779.2
$ws.Range("I76").Value = 3779.2
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3779.2
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3464.2
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3779.2
$ws.Range("I79").Value = 3779.2
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3779.2
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2687.2
$ws.Range("N79").ClearContents()
$ws.Range("H98").Value = 1850.409
$ws.Range("I98").Value = 1870.2778
$ws.Range("J98").Value = 1761
$ws.Range("K98").Value = 1870.2778
$ws.Range("L98").Value = 1761
$ws.Range("M98").Value = -372.2778000000001
$ws.Range("N98").Value = -4757
$ws.Range("H106").Value = 37040520
$ws.Range("I106").Value = 66668144
$ws.Range("K106").Value = 66668144
$ws.Range("M106").Value = -66667513
$ws.Range("H121").Value = 3330.4211
$ws.Range("J121").Value = 3330.4211
$ws.Range("L121").Value = 9991.263300000001
$ws.Range("N121").Value = -13485.2633
$ws.Range("H122").Value = 1850.409
$ws.Range("I122").Value = 1870.2778
$ws.Range("J122").Value = 1761
$ws.Range("K122").Value = 5610.8334
$ws.Range("L122").Value = 5283
$ws.Range("M122").Value = -3160.8334
$ws.Range("N122").Value = -10183
$ws.Range("H131").Value = 1499.238
$ws.Range("I131").Value = 721.25
$ws.Range("J131").Value = 3988.8
$ws.Range("K131").Value = 2163.75
$ws.Range("L131").Value = 11966.4
$ws.Range("M131").Value = 2876.25
$ws.Range("N131").Value = -22046.4
$ws.Range("H138").Value = 6322.375
$ws.Range("J138").Value = 8418.290000000001
$ws.Range("L138").Value = 25254.87
$ws.Range("N138").Value = -35534.87

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1822.7778
$ws.Range("I45").Value = 1702.6666
$ws.Range("J45").Value = 2063
$ws.Range("K45").Value = 1702.6666
$ws.Range("L45").Value = 2063
$ws.Range("M45").Value = -1325.6666
$ws.Range("N45").Value = -2817
$ws.Range("H61").Value = 10133.462
$ws.Range("I61").Value = 10818.719
$ws.Range("K61").Value = 10818.719
$ws.Range("M61").Value = -10606.719
$ws.Range("H74").Value = 1300.1428
$ws.Range("I74").Value = 540.6
$ws.Range("K74").Value = 540.6
$ws.Range("M74").Value = 333.4
$ws.Range("H77").Value = 1300.1428
$ws.Range("I77").Value = 540.6
$ws.Range("K77").Value = 2703
$ws.Range("M77").Value = 1665
$ws.Range("H109").Value = 98000
$ws.Range("J109").Value = 98000
$ws.Range("L109").Value = 98000
$ws.Range("N109").Value = -100774
$ws.Range("H110").Value = 1135802.9
$ws.Range("J110").Value = 6419.5
$ws.Range("L110").Value = 6419.5
$ws.Range("N110").Value = -10509.5
$ws.Range("H112").Value = 84332.336
$ws.Range("J112").Value = 84332.336
$ws.Range("L112").Value = 84332.336
$ws.Range("N112").Value = -87286.336
$ws.Range("H136").Value = 10133.462
$ws.Range("I136").Value = 10818.719
$ws.Range("K136").Value = 32456.157
$ws.Range("M136").Value = -29906.157

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 41670000
$ws.Range("I7").Value = 41670000
$ws.Range("K7").Value = 41670000
$ws.Range("M7").Value = -41669887
$ws.Range("H20").Value = 2359.5
$ws.Range("I20").Value = 2172.4
$ws.Range("K20").Value = 2172.4
$ws.Range("M20").Value = -1925.4
$ws.Range("H105").Value = 2419.5
$ws.Range("I105").Value = 2669.4285
$ws.Range("K105").Value = 2669.4285
$ws.Range("M105").Value = -922.4285

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2833.5
$ws.Range("I132").Value = 2886.84
$ws.Range("K132").Value = 8660.52
$ws.Range("M132").Value = -6130.52
$ws.Range("H134").Value = 2894.75
$ws.Range("I134").Value = 2524.1738
$ws.Range("J134").Value = 4599.4
$ws.Range("K134").Value = 7572.5214
$ws.Range("L134").Value = 13798.2
$ws.Range("M134").Value = -5037.5214
$ws.Range("N134").Value = -18868.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 51915484
$ws.Range("I4").Value = 1400346.1
$ws.Range("J4").Value = 178203330
$ws.Range("K4").Value = 4201038.300000001
$ws.Range("L4").Value = 534609990
$ws.Range("M4").Value = -4200926.300000001
$ws.Range("N4").Value = -534610214
$ws.Range("H12").Value = 221.77272
$ws.Range("J12").Value = 199.29411
$ws.Range("L12").Value = 597.8823299999999
$ws.Range("N12").Value = -943.8823299999999
$ws.Range("H107").Value = 1376.7
$ws.Range("J107").Value = 652.6667
$ws.Range("L107").Value = 1958.0001
$ws.Range("N107").Value = -5798.0001
$ws.Range("H138").Value = 86870.414
$ws.Range("I138").Value = 103544.8
$ws.Range("J138").Value = 3498.5
$ws.Range("K138").Value = 310634.4
$ws.Range("L138").Value = 10495.5
$ws.Range("M138").Value = -305494.4
$ws.Range("N138").Value = -20775.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1127234.9
$ws.Range("I80").Value = 2081953.2
$ws.Range("J80").Value = 36128.145
$ws.Range("K80").Value = 2081953.2
$ws.Range("L80").Value = 36128.145
$ws.Range("M80").Value = -2080955.2
$ws.Range("N80").Value = -38124.145
$ws.Range("H83").Value = 1127234.9
$ws.Range("I83").Value = 2081953.2
$ws.Range("J83").Value = 36128.145
$ws.Range("K83").Value = 10409766
$ws.Range("L83").Value = 180640.725
$ws.Range("M83").Value = -10404774
$ws.Range("N83").Value = -190624.725
$ws.Range("H102").Value = 3463.4
$ws.Range("I102").Value = 2239.3044
$ws.Range("K102").Value = 2239.3044
$ws.Range("M102").Value = -617.3044

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 4404.4
$ws.Range("J9").Value = 8000
$ws.Range("L9").Value = 8000
$ws.Range("N9").Value = -8448
$ws.Range("H22").Value = 1106.0714
$ws.Range("I22").Value = 898.75
$ws.Range("J22").Value = 1189
$ws.Range("K22").Value = 898.75
$ws.Range("L22").Value = 1189
$ws.Range("M22").Value = -603.75
$ws.Range("N22").Value = -1779
$ws.Range("H27").Value = 1106.0714
$ws.Range("I27").Value = 898.75
$ws.Range("J27").Value = 1189
$ws.Range("K27").Value = 898.75
$ws.Range("L27").Value = 1189
$ws.Range("M27").Value = -791.75
$ws.Range("N27").Value = -1403
$ws.Range("H40").Value = 4949.375
$ws.Range("I40").Value = 3336.8
$ws.Range("K40").Value = 3336.8
$ws.Range("M40").Value = -3200.8
$ws.Range("H68").Value = 670610.7
$ws.Range("I68").Value = 1264221.4
$ws.Range("J68").Value = 2798.75
$ws.Range("K68").Value = 1264221.4
$ws.Range("L68").Value = 2798.75
$ws.Range("M68").Value = -1263472.4
$ws.Range("N68").Value = -4296.75
$ws.Range("H71").Value = 670610.7
$ws.Range("I71").Value = 1264221.4
$ws.Range("J71").Value = 2798.75
$ws.Range("K71").Value = 6321107
$ws.Range("L71").Value = 13993.75
$ws.Range("M71").Value = -6317363
$ws.Range("N71").Value = -21481.75
$ws.Range("H93").Value = 1719.6666
$ws.Range("I93").Value = 1995.4
$ws.Range("K93").Value = 1995.4
$ws.Range("M93").Value = -747.4000000000001
$ws.Range("H122").Value = 6764
$ws.Range("I122").Value = 3505.3572
$ws.Range("J122").Value = 11833
$ws.Range("K122").Value = 10516.0716
$ws.Range("L122").Value = 35499
$ws.Range("M122").Value = -8066.071599999999
$ws.Range("N122").Value = -40399
$ws.Range("H132").Value = 6592
$ws.Range("I132").Value = 4904
$ws.Range("K132").Value = 14712
$ws.Range("M132").Value = -12182
$ws.Range("H136").Value = 6750.5
$ws.Range("I136").Value = 6000.8
$ws.Range("K136").Value = 18002.4
$ws.Range("M136").Value = -15452.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 100000
$ws.Range("J64").Value = 100000
$ws.Range("L64").Value = 100000
$ws.Range("N64").Value = -100496
$ws.Range("H67").Value = 100000
$ws.Range("J67").Value = 100000
$ws.Range("L67").Value = 100000
$ws.Range("N67").Value = -101716
$ws.Range("H113").Value = 740.4
$ws.Range("I113").Value = 617.7692
$ws.Range("K113").Value = 1853.3076
$ws.Range("M113").Value = 316.6924000000001
$ws.Range("H122").Value = 2320.6094
$ws.Range("I122").Value = 2278.8113
$ws.Range("K122").Value = 6836.4339
$ws.Range("M122").Value = -4386.4339
$ws.Range("H126").Value = 1456.2222
$ws.Range("I126").Value = 1308
$ws.Range("J126").Value = 1604.4445
$ws.Range("K126").Value = 3924
$ws.Range("L126").Value = 4813.333500000001
$ws.Range("M126").Value = -1454
$ws.Range("N126").Value = -9753.333500000001
